$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 530
$ws.Range("I2").Value = 530
$ws.Range("K2").Value = 530
$ws.Range("M2").Value = -417
$ws.Range("H8").Value = 351.66666
$ws.Range("I8").Value = 145.625
$ws.Range("K8").Value = 436.875
$ws.Range("M8").Value = -297.875
$ws.Range("H12").Value = 473.25
$ws.Range("I12").Value = 464.33334
$ws.Range("K12").Value = 464.33334
$ws.Range("M12").Value = -294.33334
$ws.Range("H111").Value = 33692.3
$ws.Range("I111").Value = 838.1667
$ws.Range("J111").Value = 82973.5
$ws.Range("K111").Value = 2514.5001
$ws.Range("L111").Value = 248920.5
$ws.Range("M111").Value = 552.4998999999998
$ws.Range("N111").Value = -255054.5
$ws.Range("H132").Value = 18184176
$ws.Range("I132").Value = 20409964
$ws.Range("J132").Value = 6914.3335
$ws.Range("K132").Value = 61229892
$ws.Range("L132").Value = 20743.0005
$ws.Range("M132").Value = -61227362
$ws.Range("N132").Value = -25803.0005
$ws.Range("H137").Value = 3339.4443
$ws.Range("I137").Value = 2845.4736
$ws.Range("J137").Value = 4512.625
$ws.Range("K137").Value = 8536.4208
$ws.Range("L137").Value = 13537.875
$ws.Range("M137").Value = -5986.4208
$ws.Range("N137").Value = -18637.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8134.036
$ws.Range("I32").Value = 4625.1646
$ws.Range("J32").Value = 33748.8
$ws.Range("K32").Value = 4625.1646
$ws.Range("L32").Value = 33748.8
$ws.Range("M32").Value = -4338.1646
$ws.Range("N32").Value = -34322.8
$ws.Range("H37").Value = 40999.332
$ws.Range("J37").Value = 49999.5
$ws.Range("L37").Value = 49999.5
$ws.Range("N37").Value = -50545.5
$ws.Range("H45").Value = 1532.7273
$ws.Range("J45").Value = 1846.3334
$ws.Range("L45").Value = 1846.3334
$ws.Range("N45").Value = -2600.3334
$ws.Range("H96").Value = 57999.5
$ws.Range("J96").Value = 57999.5
$ws.Range("L96").Value = 57999.5
$ws.Range("N96").Value = -63491.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2007.2727
$ws.Range("I20").Value = 2169
$ws.Range("J20").Value = 1724.25
$ws.Range("K20").Value = 2169
$ws.Range("L20").Value = 1724.25
$ws.Range("M20").Value = -1922
$ws.Range("N20").Value = -2218.25
$ws.Range("H58").Value = 25000
$ws.Range("H94").Value = 20837532
$ws.Range("H134").Value = 2452.7576
$ws.Range("I134").Value = 1727.8214
$ws.Range("J134").Value = 6512.4
$ws.Range("K134").Value = 5183.4642
$ws.Range("L134").Value = 19537.2
$ws.Range("M134").Value = -2648.4642
$ws.Range("N134").Value = -24607.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 350
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -624
$ws.Range("H7").Value = 80.5625
$ws.Range("I7").Value = 74.21429000000001
$ws.Range("K7").Value = 74.21429000000001
$ws.Range("M7").Value = 38.78570999999999
$ws.Range("H31").Value = 3250.3225
$ws.Range("I31").Value = 1853.1666
$ws.Range("K31").Value = 1853.1666
$ws.Range("M31").Value = -1558.1666
$ws.Range("H34").Value = 3250.3225
$ws.Range("I34").Value = 1853.1666
$ws.Range("K34").Value = 1853.1666
$ws.Range("M34").Value = -1651.1666
$ws.Range("H58").Value = 373382.97
$ws.Range("I58").Value = 1466.1875
$ws.Range("J58").Value = 914352.8
$ws.Range("K58").Value = 1466.1875
$ws.Range("L58").Value = 914352.8
$ws.Range("M58").Value = -1263.1875
$ws.Range("N58").Value = -914758.8
$ws.Range("H62").Value = 48371.273
$ws.Range("I62").Value = 2880.6667
$ws.Range("J62").Value = 102960
$ws.Range("K62").Value = 2880.6667
$ws.Range("L62").Value = 102960
$ws.Range("M62").Value = -2256.6667
$ws.Range("N62").Value = -104208
$ws.Range("H65").Value = 48371.273
$ws.Range("I65").Value = 2880.6667
$ws.Range("J65").Value = 102960
$ws.Range("K65").Value = 14403.3335
$ws.Range("L65").Value = 514800
$ws.Range("M65").Value = -11283.3335
$ws.Range("N65").Value = -521040
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 3032.3809
$ws.Range("I122").Value = 3547.5715
$ws.Range("K122").Value = 10642.7145
$ws.Range("M122").Value = -8192.7145
$ws.Range("H136").Value = 373382.97
$ws.Range("I136").Value = 1466.1875
$ws.Range("J136").Value = 914352.8
$ws.Range("K136").Value = 4398.5625
$ws.Range("L136").Value = 2743058.4
$ws.Range("M136").Value = -1848.5625
$ws.Range("N136").Value = -2748158.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 541.5
$ws.Range("I7").Value = 436.25
$ws.Range("K7").Value = 1308.75
$ws.Range("M7").Value = -1196.75
$ws.Range("H12").Value = 763.3333
$ws.Range("J12").Value = 915.8
$ws.Range("L12").Value = 2747.4
$ws.Range("N12").Value = -3093.4
$ws.Range("H23").Value = 1249.5555
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 1829.8334
$ws.Range("K23").Value = 267
$ws.Range("L23").Value = 5489.5002
$ws.Range("M23").Value = -32
$ws.Range("N23").Value = -5959.5002
$ws.Range("H34").Value = 1978.6471
$ws.Range("I34").Value = 170.25
$ws.Range("J34").Value = 2535.077
$ws.Range("K34").Value = 510.75
$ws.Range("L34").Value = 7605.231000000001
$ws.Range("M34").Value = -426.75
$ws.Range("N34").Value = -7773.231000000001
$ws.Range("H39").Value = 5857.3076
$ws.Range("J39").Value = 4262.0835
$ws.Range("L39").Value = 12786.2505
$ws.Range("N39").Value = -13374.2505
$ws.Range("H51").Value = 1791.8572
$ws.Range("I51").Value = 1841.8334
$ws.Range("J51").Value = 1492
$ws.Range("K51").Value = 5525.5002
$ws.Range("L51").Value = 4476
$ws.Range("M51").Value = -5065.5002
$ws.Range("N51").Value = -5396
$ws.Range("H55").Value = 2225.8125
$ws.Range("I55").Value = 299.2
$ws.Range("J55").Value = 3101.5454
$ws.Range("K55").Value = 897.5999999999999
$ws.Range("L55").Value = 9304.636200000001
$ws.Range("M55").Value = -720.5999999999999
$ws.Range("N55").Value = -9658.636200000001
$ws.Range("H70").Value = 9600
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 9600
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8577.643
$ws.Range("I70").Value = 7732.3335
$ws.Range("K70").Value = 7732.3335
$ws.Range("M70").Value = -7462.3335
$ws.Range("H73").Value = 8577.643
$ws.Range("I73").Value = 7732.3335
$ws.Range("K73").Value = 7732.3335
$ws.Range("M73").Value = -6796.3335
$ws.Range("H80").Value = 25030.375
$ws.Range("I80").Value = 8999.6
$ws.Range("J80").Value = 32317.092
$ws.Range("K80").Value = 8999.6
$ws.Range("L80").Value = 32317.092
$ws.Range("M80").Value = -8001.6
$ws.Range("N80").Value = -34313.092
$ws.Range("H83").Value = 25030.375
$ws.Range("I83").Value = 8999.6
$ws.Range("J83").Value = 32317.092
$ws.Range("K83").Value = 44998
$ws.Range("L83").Value = 161585.46
$ws.Range("M83").Value = -40006
$ws.Range("N83").Value = -171569.46
$ws.Range("H126").Value = 7443.567
$ws.Range("J126").Value = 4181.125
$ws.Range("L126").Value = 12543.375
$ws.Range("N126").Value = -17483.375
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H120").Value = 30000
$ws.Range("I120").Value = 30000
$ws.Range("K120").Value = 30000
$ws.Range("M120").Value = -25162
$ws.Range("H122").Value = 2474.9062
$ws.Range("I122").Value = 2363.3
$ws.Range("K122").Value = 7089.900000000001
$ws.Range("M122").Value = -4639.900000000001
$ws.Range("H140").Value = 106316.336
$ws.Range("J140").Value = 106316.336
$ws.Range("L140").Value = 106316.336
$ws.Range("N140").Value = -116676.336
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360
